$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Remove row 4 (old "Record 3" data row) entirely ---
$ws.Rows.Item(4).Delete() | Out-Null

# --- New header H1 (copy formatting from an existing bold/centered header cell) ---
$ws.Range("B1").Copy($ws.Range("H1")) | Out-Null
$ws.Range("H1").Value = "Number"

# --- C1 header keeps its text but gets the applyNumberFormat flag added ---
$ws.Range("C1").NumberFormat = "General"

# --- Row 2 updates ---
$ws.Range("B2").Value = "OAK"
$ws.Range("D2").Value = "Las Vegas (LAS)"
$ws.Range("F2").Value = "05/20/2020"
$ws.Range("G2").Value = "05/25/2020"
# C2 keeps the text "TRUE" value but changes number format (adds applyNumberFormat flag)
$ws.Range("C2").NumberFormat = "General"

# Move the old I2 (empty, style-only) cell content into H2 with value 4,
# keep the same number-format style (text) that I2 had.
$ws.Range("H2").NumberFormat = $ws.Range("I2").NumberFormat
$ws.Range("H2").Value = "4"
$ws.Range("I2").Clear() | Out-Null

# --- Row 3 updates ---
$ws.Range("A3").ClearContents() | Out-Null
$ws.Range("B3").Value = "SFO"
# Copy C2 (already holding text "TRUE") into C3 so the value stays textual ("TRUE"),
# then normalize its number format the same way as C2.
$ws.Range("C2").Copy($ws.Range("C3")) | Out-Null
$ws.Range("D3").Value = "Oakland (OAK)"
$ws.Range("D3").Style = "Normal"
$ws.Range("F3").Value = "05/20/2020"
$ws.Range("G3").Value = "05/25/2020"

$ws.Range("H3").NumberFormat = $ws.Range("H2").NumberFormat
$ws.Range("H3").Value = "3"

# --- Column widths: separate column E from D, give it its own width ---
$ws.Columns.Item(5).ColumnWidth = 23.5

# --- Update the active selection shown in the sheet view ---
$ws.Range("H4").Select() | Out-Null
